# feat: add 2022-Q3 data
#
# - Insert a brand-new "2022-Q3" worksheet between "总计" and "2021-Q2"
#   containing the fund holdings for that quarter.
# - Update the "总计" (summary) sheet so its first data row now reports the
#   2022-Q3 totals, and append a new row below it for the original 2021-Q2
#   totals (which used to be the only data row).

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")
$oldQtr  = $wb.Worksheets.Item("2021-Q2")

# --- 1. Create the new "2022-Q3" sheet, positioned right before "2021-Q2" ---
$oldQtr.Activate()
$newQtr = $wb.Worksheets.Add()
$newQtr.Name = "2022-Q3"

# Header row
$newQtr.Range("B1").Value = "基金代码"
$newQtr.Range("C1").Value = "基金名称"
$newQtr.Range("D1").Value = "基金规模"
$newQtr.Range("E1").Value = "股票总仓位"
$newQtr.Range("F1").Value = "仓位占比"
$newQtr.Range("G1").Value = "持有市值(亿元)"
$newQtr.Range("H1").Value = "仓位排名"

# The fund code / size / position figures are stored as TEXT (not numbers) in
# the source data, e.g. "010826" keeps its leading zero. Mark those columns
# as Text before writing so Excel doesn't auto-coerce them into numbers.
$newQtr.Range("B2:B5").NumberFormat = "@"
$newQtr.Range("D2:G5").NumberFormat = "@"

# Data rows
$newQtr.Range("A2").Value = 0
$newQtr.Range("B2").Value = "010826"
$newQtr.Range("C2").Value = "大成产业趋势混合A"
$newQtr.Range("D2").Value = "11.59"
$newQtr.Range("E2").Value = "90.99"
$newQtr.Range("F2").Value = "3.69"
$newQtr.Range("G2").Value = "0.4277"
$newQtr.Range("H2").Value = 9

$newQtr.Range("A3").Value = 1
$newQtr.Range("B3").Value = "010827"
$newQtr.Range("C3").Value = "大成产业趋势混合C"
$newQtr.Range("D3").Value = "6.29"
$newQtr.Range("E3").Value = "90.99"
$newQtr.Range("F3").Value = "3.69"
$newQtr.Range("G3").Value = "0.2321"
$newQtr.Range("H3").Value = 9

$newQtr.Range("A4").Value = 2
$newQtr.Range("B4").Value = "001735"
$newQtr.Range("C4").Value = "广发百发大数据策略成长灵活配置混合E"
$newQtr.Range("D4").Value = "0.51"
$newQtr.Range("E4").Value = "43.77"
$newQtr.Range("F4").Value = "0.59"
$newQtr.Range("G4").Value = "0.0030"
$newQtr.Range("H4").Value = 10

$newQtr.Range("A5").Value = 3
$newQtr.Range("B5").Value = "001734"
$newQtr.Range("C5").Value = "广发百发大数据策略成长灵活配置混合A"
$newQtr.Range("D5").Value = "0.27"
$newQtr.Range("E5").Value = "43.77"
$newQtr.Range("F5").Value = "0.59"
$newQtr.Range("G5").Value = "0.0016"
$newQtr.Range("H5").Value = 10

# Match the look of the "总计" sheet's header / index-column formatting
# (bold font, thin border, centered) instead of the default plain style.
$summary.Range("B1:H1").Copy()
$newQtr.Range("B1:H1").PasteSpecial(-4122)

$summary.Range("A2").Copy()
$newQtr.Range("A2:A5").PasteSpecial(-4122)

# --- 2. Update the "总计" sheet: row 2 now holds the 2022-Q3 totals ---
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.66

# --- 3. Append a new row 3 with the original 2021-Q2 totals ---
$summary.Range("A2").Copy()
$summary.Range("A3").PasteSpecial(-4122)

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2021-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0

$wb.Application.CutCopyMode = $false

# Restore "总计" as the active sheet (it was the active tab before the edit).
$summary.Activate()
